# Weekly update: insert a new price record for the Feria Lagunitas de Puerto
# Montt / Mango series. A new row is inserted at row 194 (pushing the
# existing rows 194-208 down to 195-209), and the new row is populated with
# the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 194, shifting rows 194:208 down
# to 195:209.
$ws.Rows("194:194").Insert()

# Populate the newly inserted row 194 with the new weekly observation.
$ws.Range("A194").Value = 4
$ws.Range("B194").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C194").Value = "Los Lagos"
$ws.Range("D194").Value = 44714
$ws.Range("E194").Value = 10
$ws.Range("F194").Value = "Fruta"
$ws.Range("G194").Value = 100108
$ws.Range("H194").Value = "Tropicales y subtropicales"
$ws.Range("I194").Value = 100108002
$ws.Range("J194").Value = "Mango"
$ws.Range("K194").Value = "Sin especificar"
$ws.Range("L194").Value = "Primera"
$ws.Range("M194").Value = 120
$ws.Range("N194").Value = 7500
$ws.Range("O194").Value = 8000
$ws.Range("P194").Value = 7750
$ws.Range("Q194").Value = "`$/bandeja 4 kilos"
$ws.Range("R194").Value = "Perú"
$ws.Range("S194").Value = 1938
$ws.Range("T194").Value = 4
